# Auto-generated PowerShell COM script to apply the GSW box score diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add new header cell Y1 = "PO", copying the style (bold/border/center) from X1
$ws.Range("X1").Copy() | Out-Null
$ws.Range("Y1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 25).Value2 = "PO"

# 2. For every existing data row (2-129), add an empty Y cell (text type, empty string)
for ($r = 2; $r -le 129; $r++) {
    $ws.Cells.Item($r, 25).Value2 = "'"
    $ws.Cells.Item($r, 25).Style = "Normal"
}

# 3. Give the new A130:A135 index cells the same style as the rest of column A (bold/border/center)
$ws.Range("A129").Copy() | Out-Null
$ws.Range("A130:A135").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 4. Append the six new box-score rows (130-135)
# Row 130
$ws.Cells.Item(130, 1).Value2 = 128
$ws.Cells.Item(130, 2).Value2 = "POR"
$ws.Cells.Item(130, 3).Value2 = "away"
$ws.Cells.Item(130, 4).NumberFormat = "@"
$ws.Cells.Item(130, 4).Value2 = "2025-03-10"
$ws.Cells.Item(130, 4).Style = "Normal"
$ws.Cells.Item(130, 5).Value2 = "240:00"
$ws.Cells.Item(130, 6).Value2 = 44
$ws.Cells.Item(130, 7).Value2 = 91
$ws.Cells.Item(130, 8).Value2 = 0.484
$ws.Cells.Item(130, 9).Value2 = 18
$ws.Cells.Item(130, 10).Value2 = 42
$ws.Cells.Item(130, 11).Value2 = 0.429
$ws.Cells.Item(130, 12).Value2 = 14
$ws.Cells.Item(130, 13).Value2 = 17
$ws.Cells.Item(130, 14).Value2 = 0.824
$ws.Cells.Item(130, 15).Value2 = 15
$ws.Cells.Item(130, 16).Value2 = 28
$ws.Cells.Item(130, 17).Value2 = 43
$ws.Cells.Item(130, 18).Value2 = 31
$ws.Cells.Item(130, 19).Value2 = 9
$ws.Cells.Item(130, 20).Value2 = 5
$ws.Cells.Item(130, 21).Value2 = 19
$ws.Cells.Item(130, 22).Value2 = "'"
$ws.Cells.Item(130, 22).Style = "Normal"
$ws.Cells.Item(130, 23).Value2 = 120
$ws.Cells.Item(130, 24).Value2 = -10
$ws.Cells.Item(130, 25).Value2 = 13

# Row 131
$ws.Cells.Item(131, 1).Value2 = 129
$ws.Cells.Item(131, 2).Value2 = "GSW"
$ws.Cells.Item(131, 3).Value2 = "home"
$ws.Cells.Item(131, 4).NumberFormat = "@"
$ws.Cells.Item(131, 4).Value2 = "2025-03-10"
$ws.Cells.Item(131, 4).Style = "Normal"
$ws.Cells.Item(131, 5).Value2 = "240:00"
$ws.Cells.Item(131, 6).Value2 = 45
$ws.Cells.Item(131, 7).Value2 = 85
$ws.Cells.Item(131, 8).Value2 = 0.529
$ws.Cells.Item(131, 9).Value2 = 21
$ws.Cells.Item(131, 10).Value2 = 41
$ws.Cells.Item(131, 11).Value2 = 0.512
$ws.Cells.Item(131, 12).Value2 = 19
$ws.Cells.Item(131, 13).Value2 = 21
$ws.Cells.Item(131, 14).Value2 = 0.905
$ws.Cells.Item(131, 15).Value2 = 12
$ws.Cells.Item(131, 16).Value2 = 28
$ws.Cells.Item(131, 17).Value2 = 40
$ws.Cells.Item(131, 18).Value2 = 35
$ws.Cells.Item(131, 19).Value2 = 16
$ws.Cells.Item(131, 20).Value2 = 5
$ws.Cells.Item(131, 21).Value2 = 17
$ws.Cells.Item(131, 22).Value2 = "'"
$ws.Cells.Item(131, 22).Style = "Normal"
$ws.Cells.Item(131, 23).Value2 = 130
$ws.Cells.Item(131, 24).Value2 = 10
$ws.Cells.Item(131, 25).Value2 = 16

# Row 132
$ws.Cells.Item(132, 1).Value2 = 130
$ws.Cells.Item(132, 2).Value2 = "SAC"
$ws.Cells.Item(132, 3).Value2 = "away"
$ws.Cells.Item(132, 4).NumberFormat = "@"
$ws.Cells.Item(132, 4).Value2 = "2025-03-13"
$ws.Cells.Item(132, 4).Style = "Normal"
$ws.Cells.Item(132, 5).Value2 = "240:00"
$ws.Cells.Item(132, 6).Value2 = 37
$ws.Cells.Item(132, 7).Value2 = 79
$ws.Cells.Item(132, 8).Value2 = 0.468
$ws.Cells.Item(132, 9).Value2 = 14
$ws.Cells.Item(132, 10).Value2 = 37
$ws.Cells.Item(132, 11).Value2 = 0.378
$ws.Cells.Item(132, 12).Value2 = 16
$ws.Cells.Item(132, 13).Value2 = 22
$ws.Cells.Item(132, 14).Value2 = 0.727
$ws.Cells.Item(132, 15).Value2 = 5
$ws.Cells.Item(132, 16).Value2 = 25
$ws.Cells.Item(132, 17).Value2 = 30
$ws.Cells.Item(132, 18).Value2 = 25
$ws.Cells.Item(132, 19).Value2 = 7
$ws.Cells.Item(132, 20).Value2 = 4
$ws.Cells.Item(132, 21).Value2 = 14
$ws.Cells.Item(132, 22).Value2 = 15
$ws.Cells.Item(132, 23).Value2 = 104
$ws.Cells.Item(132, 24).Value2 = -26
$ws.Cells.Item(132, 25).Value2 = "'"
$ws.Cells.Item(132, 25).Style = "Normal"

# Row 133
$ws.Cells.Item(133, 1).Value2 = 131
$ws.Cells.Item(133, 2).Value2 = "GSW"
$ws.Cells.Item(133, 3).Value2 = "home"
$ws.Cells.Item(133, 4).NumberFormat = "@"
$ws.Cells.Item(133, 4).Value2 = "2025-03-13"
$ws.Cells.Item(133, 4).Style = "Normal"
$ws.Cells.Item(133, 5).Value2 = "240:00"
$ws.Cells.Item(133, 6).Value2 = 45
$ws.Cells.Item(133, 7).Value2 = 78
$ws.Cells.Item(133, 8).Value2 = 0.577
$ws.Cells.Item(133, 9).Value2 = 22
$ws.Cells.Item(133, 10).Value2 = 39
$ws.Cells.Item(133, 11).Value2 = 0.564
$ws.Cells.Item(133, 12).Value2 = 18
$ws.Cells.Item(133, 13).Value2 = 24
$ws.Cells.Item(133, 14).Value2 = 0.75
$ws.Cells.Item(133, 15).Value2 = 7
$ws.Cells.Item(133, 16).Value2 = 36
$ws.Cells.Item(133, 17).Value2 = 43
$ws.Cells.Item(133, 18).Value2 = 35
$ws.Cells.Item(133, 19).Value2 = 8
$ws.Cells.Item(133, 20).Value2 = 5
$ws.Cells.Item(133, 21).Value2 = 14
$ws.Cells.Item(133, 22).Value2 = 21
$ws.Cells.Item(133, 23).Value2 = 130
$ws.Cells.Item(133, 24).Value2 = 26
$ws.Cells.Item(133, 25).Value2 = "'"
$ws.Cells.Item(133, 25).Style = "Normal"

# Row 134
$ws.Cells.Item(134, 1).Value2 = 132
$ws.Cells.Item(134, 2).Value2 = "NYK"
$ws.Cells.Item(134, 3).Value2 = "away"
$ws.Cells.Item(134, 4).NumberFormat = "@"
$ws.Cells.Item(134, 4).Value2 = "2025-03-15"
$ws.Cells.Item(134, 4).Style = "Normal"
$ws.Cells.Item(134, 5).Value2 = "240:00"
$ws.Cells.Item(134, 6).Value2 = 37
$ws.Cells.Item(134, 7).Value2 = 86
$ws.Cells.Item(134, 8).Value2 = 0.43
$ws.Cells.Item(134, 9).Value2 = 9
$ws.Cells.Item(134, 10).Value2 = 24
$ws.Cells.Item(134, 11).Value2 = 0.375
$ws.Cells.Item(134, 12).Value2 = 11
$ws.Cells.Item(134, 13).Value2 = 13
$ws.Cells.Item(134, 14).Value2 = 0.846
$ws.Cells.Item(134, 15).Value2 = 13
$ws.Cells.Item(134, 16).Value2 = 36
$ws.Cells.Item(134, 17).Value2 = 49
$ws.Cells.Item(134, 18).Value2 = 20
$ws.Cells.Item(134, 19).Value2 = 6
$ws.Cells.Item(134, 20).Value2 = 6
$ws.Cells.Item(134, 21).Value2 = 13
$ws.Cells.Item(134, 22).Value2 = 16
$ws.Cells.Item(134, 23).Value2 = 94
$ws.Cells.Item(134, 24).Value2 = -3
$ws.Cells.Item(134, 25).Value2 = "'"
$ws.Cells.Item(134, 25).Style = "Normal"

# Row 135
$ws.Cells.Item(135, 1).Value2 = 133
$ws.Cells.Item(135, 2).Value2 = "GSW"
$ws.Cells.Item(135, 3).Value2 = "home"
$ws.Cells.Item(135, 4).NumberFormat = "@"
$ws.Cells.Item(135, 4).Value2 = "2025-03-15"
$ws.Cells.Item(135, 4).Style = "Normal"
$ws.Cells.Item(135, 5).Value2 = "240:00"
$ws.Cells.Item(135, 6).Value2 = 33
$ws.Cells.Item(135, 7).Value2 = 83
$ws.Cells.Item(135, 8).Value2 = 0.398
$ws.Cells.Item(135, 9).Value2 = 13
$ws.Cells.Item(135, 10).Value2 = 41
$ws.Cells.Item(135, 11).Value2 = 0.317
$ws.Cells.Item(135, 12).Value2 = 18
$ws.Cells.Item(135, 13).Value2 = 22
$ws.Cells.Item(135, 14).Value2 = 0.818
$ws.Cells.Item(135, 15).Value2 = 10
$ws.Cells.Item(135, 16).Value2 = 35
$ws.Cells.Item(135, 17).Value2 = 45
$ws.Cells.Item(135, 18).Value2 = 25
$ws.Cells.Item(135, 19).Value2 = 7
$ws.Cells.Item(135, 20).Value2 = 5
$ws.Cells.Item(135, 21).Value2 = 9
$ws.Cells.Item(135, 22).Value2 = 13
$ws.Cells.Item(135, 23).Value2 = 97
$ws.Cells.Item(135, 24).Value2 = 3
$ws.Cells.Item(135, 25).Value2 = "'"
$ws.Cells.Item(135, 25).Style = "Normal"

